# edit.ps1 - apply the diff:
#   1. Rename the "AddressBook" class box on slide 1 to "TaskManager".
#   2. Refresh the cached "datetimeFigureOut" footer-date text from
#      3/17/2017 to 3/20/2017 everywhere it appears (slide master, every
#      slide layout, and the notes master).

$p = $ppt.ActivePresentation

# --- 1. AddressBook -> TaskManager -----------------------------------
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "AddressBook") {
            $sh.TextFrame.TextRange.Text = "TaskManager"
        }
    }
}

# --- helper: update any shape on a container whose text is the old date
function Update-DateShapes($shapes, $oldText, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldText) {
                $sh.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

$oldDate = "3/17/2017"
$newDate = "3/20/2017"

# --- 2a. Slide master footer date -------------------------------------
$master = $p.SlideMaster
Update-DateShapes $master.Shapes $oldDate $newDate

# --- 2b. Every slide layout footer date --------------------------------
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DateShapes $layout.Shapes $oldDate $newDate
}

# --- 2c. Notes master footer date --------------------------------------
$notesMaster = $p.NotesMaster
Update-DateShapes $notesMaster.Shapes $oldDate $newDate
